# Rename the original sheet and add five more sheets, each a copy of the
# original's layout/data, matching the commit's "Add files via upload" of a
# multi-sample workbook (WT / S47D replicate sheets).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "SRR8994357_WT"

$names = @("SRR8994358_WT", "SRR8994359_WT", "SRR8994378_S47D", "SRR8994379_S47D", "SRR8994380_S47D")
$selections = @("D15", "F7", "F9", "E12", "E11")

$prev = $ws1
for ($i = 0; $i -lt $names.Length; $i++) {
    $newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $prev)
    $newSheet.Name = $names[$i]
    $prev.Range("A1:C3").Copy($newSheet.Range("A1")) | Out-Null
    $newSheet.Columns.Item(1).ColumnWidth = $ws1.Columns.Item(1).ColumnWidth
    $newSheet.Columns.Item(2).ColumnWidth = $ws1.Columns.Item(2).ColumnWidth
    $newSheet.Columns.Item(3).ColumnWidth = $ws1.Columns.Item(3).ColumnWidth
    $newSheet.Range($selections[$i]).Select() | Out-Null
    $prev = $newSheet
}

# Restore the first sheet's own selection/active state.
$ws1.Range("D11").Select() | Out-Null
$ws1.Activate()
